$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 4580.273  # ALC H86: was 66670290
$ws.Cells.Item(86, 9).Value = 7501.5  # ALC I86: was 6001.5
$ws.Cells.Item(86, 10).Value = 3931.111  # ALC J86: was 76926330
$ws.Cells.Item(86, 11).Value = 7501.5  # ALC K86: was 6001.5
$ws.Cells.Item(86, 12).Value = 3931.111  # ALC L86: was 76926330
$ws.Cells.Item(86, 13).Value = -6378.5  # ALC M86: was -4878.5
$ws.Cells.Item(86, 14).Value = -6177.111  # ALC N86: was -76928576

$ws.Cells.Item(89, 8).Value = 4580.273  # ALC H89: was 66670290
$ws.Cells.Item(89, 9).Value = 7501.5  # ALC I89: was 6001.5
$ws.Cells.Item(89, 10).Value = 3931.111  # ALC J89: was 76926330
$ws.Cells.Item(89, 11).Value = 37507.5  # ALC K89: was 30007.5
$ws.Cells.Item(89, 12).Value = 19655.555  # ALC L89: was 384631650
$ws.Cells.Item(89, 13).Value = -31891.5  # ALC M89: was -24391.5
$ws.Cells.Item(89, 14).Value = -30887.555  # ALC N89: was -384642882

$ws.Cells.Item(96, 8).Value = 1087.2  # ALC H96: was 7822.625
$ws.Cells.Item(96, 9).Value = 1087.2  # ALC I96: was 7511.5713
$ws.Cells.Item(96, 10).Value = 0  # ALC J96: was 10000
$ws.Cells.Item(96, 11).Value = 3261.6  # ALC K96: was 22534.7139
$ws.Cells.Item(96, 12).Value = 0  # ALC L96: was 30000
$ws.Cells.Item(96, 13).ClearContents()  # ALC M96: was -21161.7139
$ws.Cells.Item(96, 14).Value = -1888.6  # ALC N96: was -32746

$ws.Cells.Item(103, 8).Value = 1181.6923  # ALC H103: was 1078.5264
$ws.Cells.Item(103, 9).Value = 1328  # ALC I103: was 1243.8182
$ws.Cells.Item(103, 10).Value = 852.5  # ALC J103: was 851.25
$ws.Cells.Item(103, 11).Value = 3984  # ALC K103: was 3731.4546
$ws.Cells.Item(103, 12).Value = 2557.5  # ALC L103: was 2553.75
$ws.Cells.Item(103, 13).Value = -3398  # ALC M103: was -3145.4546
$ws.Cells.Item(103, 14).Value = -3729.5  # ALC N103: was -3725.75

$ws.Cells.Item(110, 8).Value = 15500  # ALC H110: was 0
$ws.Cells.Item(110, 9).Value = 0  # ALC I110: was 0
$ws.Cells.Item(110, 10).Value = 15500  # ALC J110: was 0
$ws.Cells.Item(110, 11).Value = 0  # ALC K110: was 0
$ws.Cells.Item(110, 12).Value = 15500  # ALC L110: was 0
$ws.Cells.Item(110, 14).Value = -23680  # ALC N110: was None

$ws.Cells.Item(137, 8).Value = 958.9524  # ALC H137: was 959.3333
$ws.Cells.Item(137, 9).Value = 897.5294  # ALC I137: was 904.75
$ws.Cells.Item(137, 10).Value = 1220  # ALC J137: was 1134
$ws.Cells.Item(137, 11).Value = 2692.5882  # ALC K137: was 2714.25
$ws.Cells.Item(137, 12).Value = 3660  # ALC L137: was 3402
$ws.Cells.Item(137, 13).Value = -142.5882000000001  # ALC M137: was -164.25
$ws.Cells.Item(137, 14).Value = -8760  # ALC N137: was -8502

$ws.Cells.Item(141, 8).Value = 3559.5  # ALC H141: was 2638
$ws.Cells.Item(141, 9).Value = 826.4286  # ALC I141: was 809.375
$ws.Cells.Item(141, 10).Value = 9936.666999999999  # ALC J141: was 9952.5
$ws.Cells.Item(141, 11).Value = 2479.2858  # ALC K141: was 2428.125
$ws.Cells.Item(141, 12).Value = 29810.001  # ALC L141: was 29857.5
$ws.Cells.Item(141, 13).Value = 2700.7142  # ALC M141: was 2751.875
$ws.Cells.Item(141, 14).Value = -40170.001  # ALC N141: was -40217.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4906453  # ARM H32: was 4906486
$ws.Cells.Item(32, 9).Value = 5380355  # ARM I32: was 5653919.5
$ws.Cells.Item(32, 10).Value = 9466.666999999999  # ARM J32: was 6644.4443
$ws.Cells.Item(32, 11).Value = 5380355  # ARM K32: was 5653919.5
$ws.Cells.Item(32, 12).Value = 9466.666999999999  # ARM L32: was 6644.4443
$ws.Cells.Item(32, 13).Value = -5380068  # ARM M32: was -5653632.5
$ws.Cells.Item(32, 14).Value = -10040.667  # ARM N32: was -7218.4443

$ws.Cells.Item(45, 8).Value = 2210  # ARM H45: was 2226.2727
$ws.Cells.Item(45, 9).Value = 1926.25  # ARM I45: was 1948.625
$ws.Cells.Item(45, 10).Value = 2966.6667  # ARM J45: was 2966.6667
$ws.Cells.Item(45, 11).Value = 1926.25  # ARM K45: was 1948.625
$ws.Cells.Item(45, 12).Value = 2966.6667  # ARM L45: was 2966.6667
$ws.Cells.Item(45, 13).Value = -1549.25  # ARM M45: was -1571.625
$ws.Cells.Item(45, 14).Value = -3720.6667  # ARM N45: was -3720.6667

$ws.Cells.Item(122, 8).Value = 1069.6666  # ARM H122: was 910.4545000000001
$ws.Cells.Item(122, 9).Value = 1069.6666  # ARM I122: was 889.625
$ws.Cells.Item(122, 10).Value = 0  # ARM J122: was 966
$ws.Cells.Item(122, 11).Value = 3208.9998  # ARM K122: was 2668.875
$ws.Cells.Item(122, 12).Value = 0  # ARM L122: was 2898
$ws.Cells.Item(122, 13).ClearContents()  # ARM M122: was -218.875
$ws.Cells.Item(122, 14).Value = -758.9998000000001  # ARM N122: was -7798

$ws.Cells.Item(123, 8).Value = 47796.332  # ARM H123: was 33179
$ws.Cells.Item(123, 9).Value = 0  # ARM I123: was 0
$ws.Cells.Item(123, 10).Value = 47796.332  # ARM J123: was 33179
$ws.Cells.Item(123, 11).Value = 0  # ARM K123: was 0
$ws.Cells.Item(123, 12).Value = 47796.332  # ARM L123: was 33179
$ws.Cells.Item(123, 14).Value = -57596.332  # ARM N123: was -42979

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 25016174  # BSM H20: was 33354726
$ws.Cells.Item(20, 9).Value = 28303.092  # BSM I20: was 38747.312
$ws.Cells.Item(20, 10).Value = 55556904  # BSM J20: was 71430130
$ws.Cells.Item(20, 11).Value = 28303.092  # BSM K20: was 38747.312
$ws.Cells.Item(20, 12).Value = 55556904  # BSM L20: was 71430130
$ws.Cells.Item(20, 13).Value = -28056.092  # BSM M20: was -38500.312
$ws.Cells.Item(20, 14).Value = -55557398  # BSM N20: was -71430624

$ws.Cells.Item(134, 8).Value = 46168.777  # BSM H134: was 46171.223
$ws.Cells.Item(134, 9).Value = 1710.1052  # BSM I134: was 1689.5897
$ws.Cells.Item(134, 10).Value = 287515.84  # BSM J134: was 335301.84
$ws.Cells.Item(134, 11).Value = 5130.3156  # BSM K134: was 5068.7691
$ws.Cells.Item(134, 12).Value = 862547.52  # BSM L134: was 1005905.52
$ws.Cells.Item(134, 13).Value = -2595.3156  # BSM M134: was -2533.7691
$ws.Cells.Item(134, 14).Value = -867617.52  # BSM N134: was -1010975.52

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2202.2917  # CRP H31: was 2412.0466
$ws.Cells.Item(31, 9).Value = 2310.025  # CRP I31: was 2547.4722
$ws.Cells.Item(31, 10).Value = 1663.625  # CRP J31: was 1715.5714
$ws.Cells.Item(31, 11).Value = 2310.025  # CRP K31: was 2547.4722
$ws.Cells.Item(31, 12).Value = 1663.625  # CRP L31: was 1715.5714
$ws.Cells.Item(31, 13).Value = -2015.025  # CRP M31: was -2252.4722
$ws.Cells.Item(31, 14).Value = -2253.625  # CRP N31: was -2305.5714

$ws.Cells.Item(34, 8).Value = 2202.2917  # CRP H34: was 2412.0466
$ws.Cells.Item(34, 9).Value = 2310.025  # CRP I34: was 2547.4722
$ws.Cells.Item(34, 10).Value = 1663.625  # CRP J34: was 1715.5714
$ws.Cells.Item(34, 11).Value = 2310.025  # CRP K34: was 2547.4722
$ws.Cells.Item(34, 12).Value = 1663.625  # CRP L34: was 1715.5714
$ws.Cells.Item(34, 13).Value = -2108.025  # CRP M34: was -2345.4722
$ws.Cells.Item(34, 14).Value = -2067.625  # CRP N34: was -2119.5714

$ws.Cells.Item(132, 8).Value = 2315.5715  # CRP H132: was 2708.182
$ws.Cells.Item(132, 9).Value = 1838.6316  # CRP I132: was 2159.4666
$ws.Cells.Item(132, 10).Value = 3322.4443  # CRP J132: was 3884
$ws.Cells.Item(132, 11).Value = 5515.8948  # CRP K132: was 6478.399800000001
$ws.Cells.Item(132, 12).Value = 9967.332900000001  # CRP L132: was 11652
$ws.Cells.Item(132, 13).Value = -2985.8948  # CRP M132: was -3948.399800000001
$ws.Cells.Item(132, 14).Value = -15027.3329  # CRP N132: was -16712

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 910628.8  # CUL H122: was 715729.0600000001
$ws.Cells.Item(122, 9).Value = 853  # CUL I122: was 853
$ws.Cells.Item(122, 10).Value = 1430500.8  # CUL J122: was 1001679.5
$ws.Cells.Item(122, 11).Value = 7677  # CUL K122: was 7677
$ws.Cells.Item(122, 12).Value = 12874507.2  # CUL L122: was 9015115.5
$ws.Cells.Item(122, 13).Value = -5227  # CUL M122: was -5227
$ws.Cells.Item(122, 14).Value = -12879407.2  # CUL N122: was -9020015.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2551.8333  # GSM H97: was 3464
$ws.Cells.Item(97, 9).Value = 3400  # GSM I97: was 4077.7778
$ws.Cells.Item(97, 10).Value = 1946  # GSM J97: was 2359.2
$ws.Cells.Item(97, 11).Value = 3400  # GSM K97: was 4077.7778
$ws.Cells.Item(97, 12).Value = 1946  # GSM L97: was 2359.2
$ws.Cells.Item(97, 13).Value = -2904  # GSM M97: was -3581.7778
$ws.Cells.Item(97, 14).Value = -2938  # GSM N97: was -3351.2

$ws.Cells.Item(102, 8).Value = 1624.2821  # GSM H102: was 1648.6052
$ws.Cells.Item(102, 9).Value = 1211.963  # GSM I102: was 1211.963
$ws.Cells.Item(102, 10).Value = 2552  # GSM J102: was 2720.3635
$ws.Cells.Item(102, 11).Value = 1211.963  # GSM K102: was 1211.963
$ws.Cells.Item(102, 12).Value = 2552  # GSM L102: was 2720.3635
$ws.Cells.Item(102, 13).Value = 410.037  # GSM M102: was 410.037
$ws.Cells.Item(102, 14).Value = -5796  # GSM N102: was -5964.363499999999

$ws.Cells.Item(122, 8).Value = 1463733.5  # GSM H122: was 4388131.5
$ws.Cells.Item(122, 9).Value = 1881371.6  # GSM I122: was 4388131.5
$ws.Cells.Item(122, 10).Value = 2000  # GSM J122: was 0
$ws.Cells.Item(122, 11).Value = 5644114.800000001  # GSM K122: was 13164394.5
$ws.Cells.Item(122, 12).Value = 6000  # GSM L122: was 0
$ws.Cells.Item(122, 13).Value = -5641664.800000001  # GSM M122: was -13161944.5
$ws.Cells.Item(122, 14).Value = -10900  # GSM N122: was None

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1652  # LTW H93: was 970.4583
$ws.Cells.Item(93, 9).Value = 1300  # LTW I93: was 940.3182
$ws.Cells.Item(93, 10).Value = 2004  # LTW J93: was 1302
$ws.Cells.Item(93, 11).Value = 1300  # LTW K93: was 940.3182
$ws.Cells.Item(93, 12).Value = 2004  # LTW L93: was 1302
$ws.Cells.Item(93, 13).Value = -52  # LTW M93: was 307.6818
$ws.Cells.Item(93, 14).Value = -4500  # LTW N93: was -3798

$ws.Cells.Item(100, 8).Value = 11112829  # LTW H100: was 13890789
$ws.Cells.Item(100, 9).Value = 13890611  # LTW I100: was 15874873
$ws.Cells.Item(100, 10).Value = 1700  # LTW J100: was 2200
$ws.Cells.Item(100, 11).Value = 13890611  # LTW K100: was 15874873
$ws.Cells.Item(100, 12).Value = 1700  # LTW L100: was 2200
$ws.Cells.Item(100, 13).Value = -13890070  # LTW M100: was -15874332
$ws.Cells.Item(100, 14).Value = -2782  # LTW N100: was -3282

$ws.Cells.Item(122, 8).Value = 7715.05  # LTW H122: was 7715.3
$ws.Cells.Item(122, 9).Value = 8953.4  # LTW I122: was 8953.4
$ws.Cells.Item(122, 10).Value = 4000  # LTW J122: was 4001
$ws.Cells.Item(122, 11).Value = 26860.2  # LTW K122: was 26860.2
$ws.Cells.Item(122, 12).Value = 12000  # LTW L122: was 12003
$ws.Cells.Item(122, 13).Value = -24410.2  # LTW M122: was -24410.2
$ws.Cells.Item(122, 14).Value = -16900  # LTW N122: was -16903

$ws.Cells.Item(132, 8).Value = 1636.0312  # LTW H132: was 1934.7407
$ws.Cells.Item(132, 9).Value = 1114.5834  # LTW I132: was 1396.1111
$ws.Cells.Item(132, 10).Value = 3200.375  # LTW J132: was 3012
$ws.Cells.Item(132, 11).Value = 3343.7502  # LTW K132: was 4188.3333
$ws.Cells.Item(132, 12).Value = 9601.125  # LTW L132: was 9036
$ws.Cells.Item(132, 13).Value = -813.7501999999999  # LTW M132: was -1658.3333
$ws.Cells.Item(132, 14).Value = -14661.125  # LTW N132: was -14096

$ws.Cells.Item(136, 8).Value = 3736.0908  # LTW H136: was 4613.643
$ws.Cells.Item(136, 9).Value = 1982.7222  # LTW I136: was 2045.4615
$ws.Cells.Item(136, 10).Value = 11626.25  # LTW J136: was 38000
$ws.Cells.Item(136, 11).Value = 5948.1666  # LTW K136: was 6136.3845
$ws.Cells.Item(136, 12).Value = 34878.75  # LTW L136: was 114000
$ws.Cells.Item(136, 13).Value = -3398.1666  # LTW M136: was -3586.3845
$ws.Cells.Item(136, 14).Value = -39978.75  # LTW N136: was -119100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 44444.273  # WVR H46: was 56665.43
$ws.Cells.Item(46, 9).Value = 0  # WVR I46: was 0
$ws.Cells.Item(46, 10).Value = 44444.273  # WVR J46: was 56665.43
$ws.Cells.Item(46, 11).Value = 0  # WVR K46: was 0
$ws.Cells.Item(46, 12).Value = 44444.273  # WVR L46: was 56665.43
$ws.Cells.Item(46, 14).Value = -44906.273  # WVR N46: was -57127.43

$ws.Cells.Item(96, 8).Value = 7000  # WVR H96: was 3512.5
$ws.Cells.Item(96, 9).Value = 7000  # WVR I96: was 3600
$ws.Cells.Item(96, 10).Value = 0  # WVR J96: was 3250
$ws.Cells.Item(96, 11).Value = 7000  # WVR K96: was 3600
$ws.Cells.Item(96, 12).Value = 0  # WVR L96: was 3250
$ws.Cells.Item(96, 13).ClearContents()  # WVR M96: was -2227
$ws.Cells.Item(96, 14).Value = -5627  # WVR N96: was -5996

$ws.Cells.Item(103, 8).Value = 20000  # WVR H103: was 0
$ws.Cells.Item(103, 9).Value = 0  # WVR I103: was 0
$ws.Cells.Item(103, 10).Value = 20000  # WVR J103: was 0
$ws.Cells.Item(103, 11).Value = 0  # WVR K103: was 0
$ws.Cells.Item(103, 12).Value = 20000  # WVR L103: was 0
$ws.Cells.Item(103, 14).Value = -22344  # WVR N103: was None

$ws.Cells.Item(122, 8).Value = 2283.3125  # WVR H122: was 1869.8636
$ws.Cells.Item(122, 9).Value = 2132.8  # WVR I122: was 1620.75
$ws.Cells.Item(122, 10).Value = 2534.1667  # WVR J122: was 2534.1667
$ws.Cells.Item(122, 11).Value = 6398.400000000001  # WVR K122: was 4862.25
$ws.Cells.Item(122, 12).Value = 7602.500100000001  # WVR L122: was 7602.500100000001
$ws.Cells.Item(122, 13).Value = -3948.400000000001  # WVR M122: was -2412.25
$ws.Cells.Item(122, 14).Value = -12502.5001  # WVR N122: was -12502.5001

$ws.Cells.Item(123, 8).Value = 49980  # WVR H123: was 20000
$ws.Cells.Item(123, 9).Value = 0  # WVR I123: was 0
$ws.Cells.Item(123, 10).Value = 49980  # WVR J123: was 20000
$ws.Cells.Item(123, 11).Value = 0  # WVR K123: was 0
$ws.Cells.Item(123, 12).Value = 49980  # WVR L123: was 20000
$ws.Cells.Item(123, 14).Value = -59780  # WVR N123: was -29800

$ws.Cells.Item(126, 8).Value = 759.13336  # WVR H126: was 844.8461
$ws.Cells.Item(126, 9).Value = 692.125  # WVR I126: was 855.5
$ws.Cells.Item(126, 10).Value = 835.7143  # WVR J126: was 835.7143
$ws.Cells.Item(126, 11).Value = 2076.375  # WVR K126: was 2566.5
$ws.Cells.Item(126, 12).Value = 2507.1429  # WVR L126: was 2507.1429
$ws.Cells.Item(126, 13).Value = 393.625  # WVR M126: was -96.5
$ws.Cells.Item(126, 14).Value = -7447.1429  # WVR N126: was -7447.1429

$ws.Cells.Item(132, 8).Value = 1834.4231  # WVR H132: was 1903.84
$ws.Cells.Item(132, 9).Value = 1044.1333  # WVR I132: was 1111.6428
$ws.Cells.Item(132, 10).Value = 2912.0908  # WVR J132: was 2912.0908
$ws.Cells.Item(132, 11).Value = 3132.3999  # WVR K132: was 3334.9284
$ws.Cells.Item(132, 12).Value = 8736.2724  # WVR L132: was 8736.2724
$ws.Cells.Item(132, 13).Value = -602.3998999999999  # WVR M132: was -804.9284000000002
$ws.Cells.Item(132, 14).Value = -13796.2724  # WVR N132: was -13796.2724

$ws.Cells.Item(134, 8).Value = 44444.273  # WVR H134: was 56665.43
$ws.Cells.Item(134, 9).Value = 0  # WVR I134: was 0
$ws.Cells.Item(134, 10).Value = 44444.273  # WVR J134: was 56665.43
$ws.Cells.Item(134, 11).Value = 0  # WVR K134: was 0
$ws.Cells.Item(134, 12).Value = 133332.819  # WVR L134: was 169996.29
$ws.Cells.Item(134, 14).Value = -138402.819  # WVR N134: was -175066.29
